$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dimension labels become measure labels (curated dimensions reprocessed as measures)
$ws.Range("G2").Value = "iaest-measure:sexo"
$ws.Range("H2").Value = "iaest-measure:residencia-provincia-nombre"
$ws.Range("I2").Value = "iaest-measure:nacionalidad-area-nombre"
$ws.Range("J2").Value = "iaest-measure:edad-grandes-grupos"

# Row 3: role changes from "dim" to "medida" for the same columns
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("J3").Value = "medida"

# Row 4: datatype changes to xsd:int for the same columns
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"

# Row 5: mapping file references removed for sexo, nacionalidad-area-nombre, edad-grandes-grupos
$ws.Range("G5").Clear()
$ws.Range("I5").Clear()
$ws.Range("J5").Clear()
